# Add new Biduf (Bidoof) dialogue rows 28-35 to the "Все персонажи" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 27 currently closes off the previous 2-row block with the "plain"
#    style (s=4/5) and no value in column A. The new data continues that
#    block, so row 27 needs to switch to the "continuation" style (s=6/7)
#    and gain an (empty) A27 cell, matching the look of rows 6/8/11/etc.
# ---------------------------------------------------------------------------
$ws.Range("A17:E17").Copy()
$ws.Range("A27:E27").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Fill in the values for the eight new rows (28-35).
# ---------------------------------------------------------------------------

# NOTE: the source text uses the literal two-character escape sequences
# \n and \' (backslash + n, backslash + apostrophe) rather than real
# newlines/quote characters - that matches how every pre-existing row in
# this sheet already stores its script text. Single-quoted PowerShell
# strings are used throughout so those backslashes are taken literally;
# a doubled '' is PowerShell's escape for an embedded single quote.

# Row 28 - SCRIPT/G01P03A/um1204.ssb (line 1 of 3)
$ws.Range("A28").Value = 'SCRIPT/G01P03A/um1204.ssb'
$ws.Range("B28").Value = 587
$ws.Range("C28").Value = ' The great [CS:N]Dusknoir[CR] gave me\npointers on exploring, by golly.'
$ws.Range("D28").Value = ' Великий [CS:N]Даскнуар[CR] дал мне\nнаставление, ей богу.'
$ws.Range("E28").Value = ' Âåìéëéê [CS:N]Äàòëîôàñ[CR] äàì íîå\nîàòóàâìåîéå, åê áïãô.'

# Row 29 - SCRIPT/G01P03A/um1209.ssb (line 2 of 3)
$ws.Range("A29").Value = 'SCRIPT/G01P03A/um1209.ssb'
$ws.Range("B29").Value = 590
$ws.Range("C29").Value = ' He\''s given me a shot of\nconfidence in myself, yup yup!'
$ws.Range("D29").Value = ' Он помог мне поверить в себя,\nда-да!'
$ws.Range("E29").Value = ' Ïî ðïíïã íîå ðïâåñéóû â òåáÿ,\näà-äà!'

# Row 30 - continuation (line 3 of 3, same block as row 29)
$ws.Range("B30").Value = 593
$ws.Range("C30").Value = ' I\''ve got the spirit for a great\nday of training! It\''s called being positive,\nby gosh!'
$ws.Range("D30").Value = ' У меня появился настрой на\nотличный день тренировок! Ей богу, вот\nэто позитивный подход к делу!'
$ws.Range("E30").Value = ' Ô íåîÿ ðïÿâéìòÿ îàòóñïê îà\nïóìéœîúê äåîû óñåîéñïâïë! Åê áïãô, âïó\nüóï ðïèéóéâîúê ðïäöïä ë äåìô!'

# Row 31 - SCRIPT/T01P01A/um1306.ssb
$ws.Range("A31").Value = 'SCRIPT/T01P01A/um1306.ssb'
$ws.Range("B31").Value = 568
$ws.Range("C31").Value = ' I have to get a move on and\npack... Huff-huff…'
$ws.Range("D31").Value = ' Мне нужно поскорее собрать вещи...\nУхх-ухх...'
$ws.Range("E31").Value = ' Íîå îôçîï ðïòëïñåå òïáñàóû âåþé...\nÔöö-ôöö…'

# Row 32 - SCRIPT/P01P01A/um1302.ssb
$ws.Range("A32").Value = 'SCRIPT/P01P01A/um1302.ssb'
$ws.Range("B32").Value = 549
$ws.Range("C32").Value = ' Um...[K] Have you all happened to\nsee [CS:N]Dugtrio[CR] anywhere?'
$ws.Range("D32").Value = ' Эм...[K] Вы, случаем, не видели [CS:N]Дагтрио[CR]?'
$ws.Range("E32").Value = ' Üí...[K] Âú, òìôœàåí, îå âéäåìé [CS:N]Äàãóñéï[CR]?'

# Row 33 - SCRIPT/G01P03A/um1510.ssb (line 1 of 2)
$ws.Range("A33").Value = 'SCRIPT/G01P03A/um1510.ssb'
$ws.Range("B33").Value = 527
$ws.Range("C33").Value = ' I\''ll get out there again as soon\nas I get geared up, yup yup!'
$ws.Range("D33").Value = ' Я отправлюсь туда сразу же как\nподготовлюсь, да-да!'
$ws.Range("E33").Value = ' Ÿ ïóðñàâìýòû óôäà òñàèô çå ëàë\nðïäãïóïâìýòû, äà-äà!'

# Row 34 - continuation (line 2 of 2, same block as row 33)
$ws.Range("B34").Value = 530
$ws.Range("C34").Value = ' Off to [CS:P]Crystal Crossing[CR]!\nYup yup!'
$ws.Range("D34").Value = ' Прямиком в [CS:P]Кристальный Переход[CR]!\nДа-да!'
$ws.Range("E34").Value = ' Ðñÿíéëïí â [CS:P]Ëñéòóàìûîúê Ðåñåöïä[CR]!\nÄà-äà!'

# Row 35 - SCRIPT/T01P01A/um1605.ssb
$ws.Range("A35").Value = 'SCRIPT/T01P01A/um1605.ssb'
$ws.Range("B35").Value = 508
$ws.Range("C35").Value = ' I\''ve been spreading that rumor\nthe best I could, yup yup!'
$ws.Range("D35").Value = ' Я всеми силами распространяю слухи, да-да!'
$ws.Range("E35").Value = ' Ÿ âòåíé òéìàíé ñàòñðïòóñàîÿý òìôöé, äà-äà!'

# ---------------------------------------------------------------------------
# 3) Apply the right cell formatting (border/font/wrap style) to each new
#    row by copying formats from an existing row that already uses the
#    matching style:
#      - rows 28, 29, 33, 35   -> "plain" style   (like row 26, s=4/5)
#      - rows 30, 34           -> "continuation"  (like row 17, s=6/7)
#      - rows 31, 32           -> "highlight"     (like row 15, s=8/9)
# ---------------------------------------------------------------------------
$ws.Range("A26:E26").Copy()
$ws.Range("A28:E28").PasteSpecial(-4122)
$ws.Range("A29:E29").PasteSpecial(-4122)
$ws.Range("A33:E33").PasteSpecial(-4122)
$ws.Range("A35:E35").PasteSpecial(-4122)

$ws.Range("A17:E17").Copy()
$ws.Range("A30:E30").PasteSpecial(-4122)
$ws.Range("A34:E34").PasteSpecial(-4122)

$ws.Range("A15:E15").Copy()
$ws.Range("A31:E31").PasteSpecial(-4122)
$ws.Range("A32:E32").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Row heights, matching the authored values.
# ---------------------------------------------------------------------------
$ws.Rows.Item(28).RowHeight = 43.2
$ws.Rows.Item(29).RowHeight = 43.2
$ws.Rows.Item(30).RowHeight = 31.8
$ws.Rows.Item(31).RowHeight = 43.2
$ws.Rows.Item(32).RowHeight = 43.2
$ws.Rows.Item(33).RowHeight = 43.2
$ws.Rows.Item(34).RowHeight = 21.6
$ws.Rows.Item(35).RowHeight = 43.2

# ---------------------------------------------------------------------------
# 5) Scroll/selection state, as recorded in the workbook after the edit.
# ---------------------------------------------------------------------------
$ws.Range("C32").Select()
$excel.ActiveWindow.ScrollRow = 31
